$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "30.648.06"
$ws.Range("E2").Value = "  +2.35%  "
$ws.Range("D3").Value = "1.675.01"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  +0.18%  "
Set-TextValue $ws.Range("D5") "219.63"
$ws.Range("E5").Value = "  +2.23%  "
Set-TextValue $ws.Range("D6") "0.528"
$ws.Range("E6").Value = "  +1.98%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.19%  "
Set-TextValue $ws.Range("D8") "29.28"
$ws.Range("E8").Value = "  +1.56%  "
Set-TextValue $ws.Range("D9") "0.264"
$ws.Range("E9").Value = "  +2.30%  "
Set-TextValue $ws.Range("D10") "0.0639"
$ws.Range("E10").Value = "  +5.00%  "
Set-TextValue $ws.Range("D11") "0.0904"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "1.916.84"
$ws.Range("E12").Value = "  +2.66%  "
$ws.Range("D13").Value = "1.680.81"
$ws.Range("E13").Value = "  +2.82%  "
Set-TextValue $ws.Range("D14") "0.606"
$ws.Range("E14").Value = "  +7.50%  "
Set-TextValue $ws.Range("D15") "9.86"
$ws.Range("E15").Value = "  +5.93%  "
Set-TextValue $ws.Range("D16") "4.02"
$ws.Range("E16").Value = "  +4.41%  "
$ws.Range("D17").Value = "30.663.83"
$ws.Range("E17").Value = "  +2.38%  "
Set-TextValue $ws.Range("D18") "66.21"
$ws.Range("E18").Value = "  +3.15%  "
Set-TextValue $ws.Range("D19") "241.72"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "0.0₃0720"
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("E21").Value = "  -0.01%  "
Set-TextValue $ws.Range("D22") "4.23"
$ws.Range("E22").Value = "  +2.34%  "
Set-TextValue $ws.Range("D23") "9.96"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("E24").Value = "  -0.62%  "
Set-TextValue $ws.Range("D25") "159.06"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("E26").Value = "  +2.87%  "
Set-TextValue $ws.Range("D27") "15.79"
$ws.Range("E27").Value = "  +1.90%  "
Set-TextValue $ws.Range("D28") "6.67"
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("E29").Value = "  +0.26%  "
Set-TextValue $ws.Range("D30") "0.0493"
$ws.Range("E30").Value = "  +0.34%  "
Set-TextValue $ws.Range("D31") "1.15"
$ws.Range("E31").Value = "  +3.51%  "
Set-TextValue $ws.Range("D32") "3.46"
$ws.Range("E32").Value = "  +1.86%  "
Set-TextValue $ws.Range("D33") "3.33"
$ws.Range("E33").Value = "  +4.76%  "
$ws.Range("D34").Value = "1.498.35"
$ws.Range("E34").Value = "  +4.72%  "
Set-TextValue $ws.Range("D35") "1.77"
$ws.Range("E35").Value = "  +6.84%  "
Set-TextValue $ws.Range("D36") "83.47"
$ws.Range("E36").Value = "  +9.92%  "
$ws.Range("E37").Value = "  -0.57%  "
Set-TextValue $ws.Range("D38") "0.596"
$ws.Range("E38").Value = "  +7.73%  "
$ws.Range("E39").Value = "  +4.32%  "
$ws.Range("E40").Value = "  -3.34%  "
Set-TextValue $ws.Range("D41") "2.30"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("E42").Value = "  +1.07%  "
Set-TextValue $ws.Range("D43") "0.837"
$ws.Range("E43").Value = "  +0.72%  "
Set-TextValue $ws.Range("D44") "0.0500"
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("E46").Value = "  +0.14%  "
Set-TextValue $ws.Range("D47") "5.55"
$ws.Range("E47").Value = "  +3.66%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.811.26"
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D49") "49.88"
$ws.Range("E49").Value = "  -2.59%  "
Set-TextValue $ws.Range("D50") "93.58"
$ws.Range("E50").Value = "  +3.49%  "
$ws.Range("E51").Value = "  +1.77%  "
